# "color slope in leachate data"
# Fill in the missing leachate measurements for rows 20-55 of the
# "Rainfall 4" sheet, highlight the newly entered G/H (slope) columns in
# yellow, and extend the I-column (H-G) formula down through row 55.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Rows 20-37 were sampled 2021-06-22, rows 38-55 were sampled 2021-06-24.
$firstDate  = 20210622
$secondDate = 20210624

for ($r = 20; $r -le 55; $r++) {
    if ($r -le 37) {
        $sampleDate = $firstDate
    } else {
        $sampleDate = $secondDate
    }

    $ws.Range("D${r}").Value = $sampleDate
    $ws.Range("J${r}").Value = $sampleDate

    # New measured values, highlighted with a yellow fill.
    $ws.Range("G${r}").Value = 100
    $ws.Range("H${r}").Value = 1000
    $ws.Range("G${r}:H${r}").Interior.Color = 65535

    # Extend the "H - G" leachate volume formula down into this row.
    $ws.Range("I${r}").Formula = "=H${r}-G${r}"
}

$excel.Calculate()

# Match the author's final selection/cursor position.
$ws.Range("J34").Select()
